# Auto-generated edit script applying the diff to the two affected sheets
$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Failed Audits-Closed Tickets")
$ws = $ws2
# Row 2
$ws.Range('G2').Value = 'T-82013'
$ws.Range('J2').Value = 'FLIK 3C - HACCP thermometer in place: Service/Line temperatures'
$ws.Range('K2').Value = 'logs filled out daily'
$ws.Range('L2').Value = 2

# Row 3
$ws.Range('G3').Value = 'T-82008'
$ws.Range('J3').Value = 'FLIK 3A - Paper Logs: Service/line temperatures'
$ws.Range('K3').Value = 'complete logs for Entree'
$ws.Range('L3').Value = 1

# Row 4
$ws.Range('G4').Value = 'T-82009'
$ws.Range('J4').Value = '812--When time is observed being used as a public health control: written procedures are available and followed.'
$ws.Range('K4').Value = 'time temp log will be used for breakfast sand'
$ws.Range('L4').Value = 1

# Row 5
$ws.Range('G5').Value = 'T-82012'
$ws.Range('J5').Value = '802--TCS food is cooled following proper methods which facilitate rapid cooling. All TCS food is properly chilled to 41F or less prior to placing into service or display.'
$ws.Range('K5').Value = 'proper cooling technique used'
$ws.Range('L5').Value = 1

# Row 6
$ws.Range('G6').Value = 'T-82010'
$ws.Range('J6').Value = 'FLIK 7 - Paper Logs: Warewashing'
$ws.Range('K6').Value = 'ware washing logs implemented'
$ws.Range('L6').Value = 1

# Row 7
$ws.Range('G7').Value = 'T-82007'
$ws.Range('J7').Value = 'FLIK 6 - Paper Logs: Cooling TCS foods are monitored/recorded'
$ws.Range('K7').Value = 'cooking cooling logs will be used'
$ws.Range('L7').Value = 1

# Row 9
$ws.Range('G9').Value = 'T-82014'
$ws.Range('J9').Value = '809--TCS foods are held cold at 41F or below - Mechanical Equipment/Facilities Related.'
$ws.Range('K9').Value = 'facilities cleaned grab and go cooler'
$ws.Range('L9').Value = 1

# Row 12
$ws.Range('G12').Value = 'T-81690'
$ws.Range('J12').Value = '1306-- Floors, walls, ceilings, and attached equipment are maintained in good condition and clean.'
$ws.Range('K12').Value = 'The holes observed (under three bay sink, by dishwasher entrance, by hand sink near cooler 1, ceiling tile above three bay sink) have had work orders placed to have them repaired/plugged'
$ws.Range('L12').Value = 0

# Row 13
$ws.Range('G13').Value = 'T-81692'
$ws.Range('J13').Value = 'FLIK 3C - HACCP thermometer in place: Service/Line temperatures'
$ws.Range('K13').Value = 'Immediately following service today there was a retraining with all staff to review the HACCP monitor SOPs, and it was stressed to hit upload after each batch of temperatures are measured. This will be echoed in preservice meetings, weekly staff meetings, and a daily check in with each staff member when they turn in the paper logs of the day and tablet is reviewed for uploaded information.'
$ws.Range('L13').Value = 0

# Row 14
$ws.Range('G14').Value = 'T-81691'
$ws.Range('J14').Value = '506--Food is protected from contamination from customers during service and display/consumer self-service.'
$ws.Range('K14').Value = 'With peanut butter being available to the entire school coupled with the severity and seriousness of the risk involved with the open container of peanut butter on the deli station, the implementation of a stand alone station for peanut butter is necessary. An email to the client will be sent to inform/explain the change to procedures and a new isolated service area will be utilized, conforming to the company standards'
$ws.Range('L14').Value = 0

# Row 16
$ws.Range('G16').Value = 'T-81694'
$ws.Range('J16').Value = 'FLIK 3B - HACCP thermometer in place: Food Production '
$ws.Range('K16').Value = 'A daily check in with each station at the end of service will be reviewing logs performed by the hourly staff by either the DDS or EC to ensure that all sheets have been submitted and are organized/on file to avoid this issue moving forward'
$ws.Range('L16').Value = 0

# Row 17
$ws.Range('G17').Value = 'T-81689'
$ws.Range('J17').Value = '1303--Plumbing is maintained in good repair.'
$ws.Range('K17').Value = 'Work order has been placed with facilities and will be repaired promptly'
$ws.Range('L17').Value = 0

$ws4 = $wb.Worksheets.Item("All Audits-Open Tickets")
$ws = $ws4
# Row 2
$ws.Range('G2').Value = 'T-81951'
$ws.Range('J2').Value = '808--TCS foods are held cold at 41F or below - Behavior/Process Related.'

# Row 3
$ws.Range('G3').Value = 'T-81945'
$ws.Range('J3').Value = '402--Food is in wholesome, sound condition and is unadulterated (no spoilage, no contamination with foreign material, intact food packaging). Damaged and/or recalled items are segregated.'

# Row 5
$ws.Range('G5').Value = 'T-81943'
$ws.Range('J5').Value = '104--The unit has the Compass Group written procedures for responding to a diarrheal or vomiting event.'

# Row 6
$ws.Range('B6').Value = "'61138"
$ws.Range('B6').ClearFormats()
$ws.Range('C6').Value = 'Castilleja School'
$ws.Range('G6').Value = 'T-81947'

# Row 7
$ws.Range('B7').Value = "'53051"
$ws.Range('B7').ClearFormats()
$ws.Range('C7').Value = 'The Logan School'
$ws.Range('G7').Value = 'T-81942'

# Row 8
$ws.Range('G8').Value = 'T-81944'

# Row 9
$ws.Range('G9').Value = 'T-81948'
$ws.Range('J9').Value = '808--TCS foods are held cold at 41F or below - Behavior/Process Related.'

# Row 10
$ws.Range('G10').Value = 'T-81950'
$ws.Range('J10').Value = '600--Food contact surfaces are properly cleaned and sanitized at regular intervals and are clean to sight and touch.'

# Row 11
$ws.Range('G11').Value = 'T-81946'
$ws.Range('J11').Value = '505--Food is properly stored to protect it from contamination (i.e. stored off the floor, kept covered during storage, properly stored on ice, not located in prohibited areas).'

# Row 14
$ws.Range('G14').Value = 'T-81411'
$ws.Range('J14').Value = '901-- Food containers holding food/ingredients that are removed from their original packages are labeled to identify contents that are not easily identifiable'

# Row 15
$ws.Range('G15').Value = 'T-81412'
$ws.Range('J15').Value = '1001--Clean utensils and equipment are properly stored, air-dried, and handled. Clean tableware is properly handled and protected.'

# Row 18
$ws.Range('G18').Value = 'T-81342'
$ws.Range('J18').Value = 'FLIK 4A - Paper Logs: TCS Receiving temperatures documented for each delivery'

# Row 19
$ws.Range('G19').Value = 'T-81332'
$ws.Range('J19').Value = '600--Food contact surfaces are properly cleaned and sanitized at regular intervals and are clean to sight and touch.'

# Row 20
$ws.Range('G20').Value = 'T-81337'
$ws.Range('J20').Value = 'FLIK 3A - Paper Logs: Service/line temperatures'

# Row 21
$ws.Range('G21').Value = 'T-81328'
$ws.Range('J21').Value = 'FLIK 5A - Paper Logs: Refrigeration temperatures are monitored AM and PM daily'

# Row 23
$ws.Range('G23').Value = 'T-81339'
$ws.Range('J23').Value = 'FLIK 3B - HACCP thermometer in place: Food Production '

# Row 24
$ws.Range('G24').Value = 'T-81343'
$ws.Range('J24').Value = '506--Food is protected from contamination from customers during service and display/consumer self-service.'

# Row 26
$ws.Range('G26').Value = 'T-81327'
$ws.Range('J26').Value = '1300--Pest activity/evidence of pest activity is not observed.'

# Row 27
$ws.Range('G27').Value = 'T-81331'
$ws.Range('J27').Value = '600--Food contact surfaces are properly cleaned and sanitized at regular intervals and are clean to sight and touch.'

# Row 28
$ws.Range('G28').Value = 'T-81340'
$ws.Range('J28').Value = 'FLIK 7 - Paper Logs: Warewashing'

# Row 29
$ws.Range('G29').Value = 'T-81333'
$ws.Range('J29').Value = 'FLIK 4B - HACCP Thermometer in place: TCS Receiving Temperatures documented for each delivery'

# Row 30
$ws.Range('G30').Value = 'T-81341'
$ws.Range('J30').Value = '100--Person-in-charge is present and is certified by an accredited program and the certificate is current.'

# Row 32
$ws.Range('G32').Value = 'T-81329'
$ws.Range('J32').Value = 'FLIK 5B - HACCP Thermometer in place: Walk in coolers or Primary Storage Coolers are monitored AM and PM daily'

# Row 33
$ws.Range('G33').Value = 'T-81338'
$ws.Range('J33').Value = 'FLIK 6 - Paper Logs: Cooling TCS foods are monitored/recorded'

# Row 34
$ws.Range('G34').Value = 'T-81330'
$ws.Range('J34').Value = 'FLIK 3A - Paper Logs: Food Production Temperatures'

